$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("sheet1")
$ws2 = $wb.Worksheets.Item("After keeping max 26 vids")

# ---------------------------------------------------------------
# sheet1: row 8 / row 9 value edits
# ---------------------------------------------------------------
$ws1.Range("C8").Value = 35

# ---------------------------------------------------------------
# sheet2: new footnote cell (written early so the shared-string
# table picks up the ids in the same order the original commit did)
# ---------------------------------------------------------------
$ws2.Range("C58").Value = "All sentences are either 23, 24, 25, 26 vids"

# back to sheet1
$ws1.Range("I8").Value = "*P2 has 2 P3 vids (extra"
$ws1.Range("C9").Value = 38
$ws1.Range("B60").Value = "P6 blooper vid: 06_0030_(08_04_21_22_19_30)_c"

# ---------------------------------------------------------------
# sheet2: Q4 edit
# ---------------------------------------------------------------
$ws2.Range("Q4").Value = 6

# ---------------------------------------------------------------
# sheet2: row 51 -- new border style (top/left/right, no bottom)
# ---------------------------------------------------------------
$src51 = $ws1.Range("C2")
$src51.Copy()
$row51 = $ws2.Range("C51:H51")
$row51.PasteSpecial(-4122)
$row51.Borders.Item(9).LineStyle = -4142
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# sheet2: row 52 -- Total row, new formulas + new cells
# ---------------------------------------------------------------
$ws2.Range("B52").Value = "Total"
$ws2.Range("C52").Formula = "=SUM(C2:C51)"
$ws2.Range("D52").Formula = "=SUM(D2:D51)"
$ws2.Range("E52").Formula = "=SUM(E2:E51)"
$ws2.Range("F52").Formula = "=SUM(F2:F51)"
$ws2.Range("G52").Formula = "=SUM(G2:G51)"
$ws2.Range("H52").Formula = "=SUM(H2:H51)"

$ws2.Range("L52").Formula = "=SUM(L2:L51)"
$ws2.Range("M52").Formula = "=SUM(M2:M51)"
$ws2.Range("N52").Formula = "=SUM(N2:N51)"
$ws2.Range("O52").Formula = "=SUM(O2:O51)"
$ws2.Range("P52").Formula = "=SUM(P2:P51)"
$ws2.Range("Q52").Formula = "=SUM(Q2:Q51)"

# I52 / R52 pick up the new right/top/bottom (no left) border style
$srcI52 = $ws1.Range("C2")
$srcI52.Copy()
$ws2.Range("I52").PasteSpecial(-4122)
$ws2.Range("I52").Borders.Item(7).LineStyle = -4142
$excel.CutCopyMode = 0
$ws2.Range("I52").Value = "Total"

$srcR52 = $ws1.Range("C2")
$srcR52.Copy()
$ws2.Range("R52").PasteSpecial(-4122)
$ws2.Range("R52").Borders.Item(7).LineStyle = -4142
$excel.CutCopyMode = 0
$ws2.Range("R52").Value = "Total"

# ---------------------------------------------------------------
# sheet2: row 53 -- new cumulative row
# ---------------------------------------------------------------
$ws2.Range("B53").Value = "Cumulative"
$ws2.Range("C53").Formula = "=SUM(C52)"
$ws2.Range("D53").Formula = "=SUM(C52:D52)"
$ws2.Range("E53").Formula = "=SUM(C52:E52)"
$ws2.Range("F53").Formula = "=SUM(C52:F52)"
$ws2.Range("G53").Formula = "=SUM(C52:G52)"
$ws2.Range("H53").Formula = "=SUM(C52:H52)"

$ws2.Range("L53").Formula = "=SUM(L52)"
$ws2.Range("M53").Formula = "=SUM(L52:M52)"
$ws2.Range("N53").Formula = "=SUM(L52:N52)"
$ws2.Range("O53").Formula = "=SUM(L52:O52)"
$ws2.Range("P53").Formula = "=SUM(L52:P52)"
$ws2.Range("Q53").Formula = "=SUM(L52:Q52)"

# ---------------------------------------------------------------
# sheet2: row 56 -- All Videos total now sums J52 and S52 only
# (L52:Q52 are no longer blank, so SUM(J52:S52) would double count)
# ---------------------------------------------------------------
$ws2.Range("J56").Formula = "=SUM(J52, S52)"

# ---------------------------------------------------------------
# Column widths on sheet2 (U / W got wider after the new content)
# ---------------------------------------------------------------
$ws2.Columns.Item(21).AutoFit()
$ws2.Columns.Item(23).AutoFit()

# ---------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.Zoom = 62
$ws1.Range("O4").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 86
$ws2.Range("L54:Q54").Select()
